# Reduce the title MERGEFIELD ("TitleLabel") font size from 28pt (sz=56,
# inherited from the "Titel" paragraph style) down to 18pt by stamping an
# explicit run-level override (w:sz=36 / w:szCs=36 half-points) on every run
# in that paragraph (the field-char runs, the instrText run and the
# field-result run), plus the paragraph mark's rPr, exactly mirroring what
# Word itself writes when you select the title and change the Font Size box
# to 18 with complex-script sizing kept in sync.
#
# Word's Range.Font.Size COM property only ever stamps w:sz (western size);
# it intentionally leaves w:szCs alone, so a plain
#     $range.Font.Size = 18
# cannot reproduce the paired w:sz/w:szCs the diff expects. Range.InsertXML
# (a real Word COM method - "InsertXML REPLACES that range's contents") lets
# us hand Word the exact finished OOXML for the paragraph instead, which is
# the most faithful way to land both attributes through the object model.

$d = $word.ActiveDocument

# --- 1) Title paragraph: MERGEFIELD TitleLabel -> shrink to 18pt (sz/szCs=36) ---
$titlePara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "TitleLabel") {
        $titlePara = $para
        break
    }
}
if ($titlePara -eq $null) {
    throw "Could not locate the TitleLabel MERGEFIELD paragraph"
}

$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="58025E90" w14:textId="4B7EE872" w:rsidR="00E36D63" w:rsidRDefault="0041346C" w:rsidP="00E36D63"><w:pPr><w:pStyle w:val="Titel"/><w:spacing w:after="240"/><w:rPr><w:noProof/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD  TitleLabel  \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00F262D6"><w:rPr><w:noProof/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>«TitleLabel»</w:t></w:r><w:r><w:rPr><w:noProof/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titlePara.Range.InsertXML($titleXml)

# --- 2) EmailLabel field run: drop the stale <w:lastRenderedPageBreak/> marker ---
# (Word repaginates and regenerates these cached page-break markers itself;
# now that the title is shorter the old break position is stale, so the
# marker is removed from the run that led with it.)
$emailPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "EmailLabel") {
        $emailPara = $para
        break
    }
}
if ($emailPara -eq $null) {
    throw "Could not locate the EmailLabel MERGEFIELD paragraph"
}

$emailXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="0BFB350E" w14:textId="578BDBC7" w:rsidR="00780652" w:rsidRPr="003B0EAD" w:rsidRDefault="00856BF5" w:rsidP="00666AC3"><w:pPr><w:rPr><w:rFonts w:ascii="HelveticaNeue LT 45 Light" w:hAnsi="HelveticaNeue LT 45 Light" w:cs="Arial"/><w:noProof/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="HelveticaNeue LT 45 Light" w:hAnsi="HelveticaNeue LT 45 Light" w:cs="Arial"/><w:noProof/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="HelveticaNeue LT 45 Light" w:hAnsi="HelveticaNeue LT 45 Light" w:cs="Arial"/><w:noProof/><w:szCs w:val="20"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD  EmailLabel  \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="HelveticaNeue LT 45 Light" w:hAnsi="HelveticaNeue LT 45 Light" w:cs="Arial"/><w:noProof/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00F262D6"><w:rPr><w:rFonts w:ascii="HelveticaNeue LT 45 Light" w:hAnsi="HelveticaNeue LT 45 Light" w:cs="Arial"/><w:noProof/><w:szCs w:val="20"/></w:rPr><w:t>«EmailLabel»</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="HelveticaNeue LT 45 Light" w:hAnsi="HelveticaNeue LT 45 Light" w:cs="Arial"/><w:noProof/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$emailPara.Range.InsertXML($emailXml)

Write-Host "Title font size reduced to 18pt; stale lastRenderedPageBreak removed."
